$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F14").Value = 1

$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 25
$ws.Range("I15").Value = 19
$ws.Range("J15").Value = 27
$ws.Range("K15").Value = -29.629629629629
$ws.Range("L15").Value = -13.636363636363
$ws.Range("M15").Value = 35.714285714285
$ws.Range("N15").Value = -57.777777777777

$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 32
$ws.Range("G16").Value = 38
$ws.Range("H16").Value = -15.78947368421
$ws.Range("I16").Value = 350
$ws.Range("J16").Value = 304
$ws.Range("K16").Value = 15.131578947368
$ws.Range("L16").Value = 9.375
$ws.Range("M16").Value = 30.597014925373
$ws.Range("N16").Value = -56.79012345679

$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 112.5
$ws.Range("F17").Value = 61
$ws.Range("G17").Value = 68
$ws.Range("H17").Value = -10.294117647058
$ws.Range("I17").Value = 549
$ws.Range("J17").Value = 521
$ws.Range("K17").Value = 5.374280230326
$ws.Range("L17").Value = 16.313559322033
$ws.Range("M17").Value = 98.91304347826
$ws.Range("N17").Value = 6.189555125725

$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -8.333333333333
$ws.Range("I18").Value = 226
$ws.Range("J18").Value = 182
$ws.Range("K18").Value = 24.175824175824
$ws.Range("L18").Value = 21.505376344086
$ws.Range("M18").Value = 26.966292134831
$ws.Range("N18").Value = -76.55601659751

$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 14.035087719298
$ws.Range("I19").Value = 463
$ws.Range("J19").Value = 380
$ws.Range("K19").Value = 21.842105263157
$ws.Range("L19").Value = 39.039039039039
$ws.Range("M19").Value = 91.322314049586
$ws.Range("N19").Value = 34.593023255814

$ws.Range("C20").Value = 7
$ws.Range("E20").Value = 16.666666666666
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 10.714285714285
$ws.Range("I20").Value = 228
$ws.Range("J20").Value = 252
$ws.Range("K20").Value = -9.523809523809
$ws.Range("L20").Value = 8.056872037914
$ws.Range("M20").Value = 171.428571428571
$ws.Range("N20").Value = -42.713567839196

$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = 10.63829787234
$ws.Range("F21").Value = 217
$ws.Range("G21").Value = 219
$ws.Range("H21").Value = -0.913242009132
$ws.Range("I21").Value = 1841
$ws.Range("J21").Value = 1678
$ws.Range("K21").Value = 9.713945172824
$ws.Range("L21").Value = 18.31619537275
$ws.Range("M21").Value = 72.863849765258
$ws.Range("N21").Value = -40.727623953638

$ws.Range("D23").Value = 2
$ws.Range("G23").Value = 4
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = 10
$ws.Range("L23").Value = -8.333333333333

$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -38.709677419354
$ws.Range("F24").Value = 104
$ws.Range("G24").Value = 78
$ws.Range("H24").Value = 33.333333333333
$ws.Range("I24").Value = 791
$ws.Range("J24").Value = 750
$ws.Range("K24").Value = 5.466666666666
$ws.Range("L24").Value = -13.928182807399
$ws.Range("M24").Value = 9.556786703601

$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 40
$ws.Range("H25").Value = 150
$ws.Range("I25").Value = 237
$ws.Range("J25").Value = 253
$ws.Range("K25").Value = -6.324110671936
$ws.Range("L25").Value = -43.301435406698

$ws.Range("C26").Value = 34
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = 30.76923076923
$ws.Range("F26").Value = 86
$ws.Range("G26").Value = 94
$ws.Range("H26").Value = -8.510638297872
$ws.Range("I26").Value = 739
$ws.Range("J26").Value = 676
$ws.Range("K26").Value = 9.319526627218
$ws.Range("L26").Value = 9.158050221565
$ws.Range("M26").Value = 0.544217687074

$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 33
$ws.Range("J27").Value = 41
$ws.Range("K27").Value = -19.512195121951
$ws.Range("L27").Value = -25

$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("J28").Value = 54
$ws.Range("K28").Value = 1.851851851851
$ws.Range("L28").Value = -6.779661016949

$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 31
$ws.Range("K29").Value = 10.714285714285
$ws.Range("L29").Value = -29.545454545454
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -67.708333333333

$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 50
$ws.Range("I30").Value = 26
$ws.Range("K30").Value = 23.809523809523
$ws.Range("L30").Value = -31.578947368421
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = -68.67469879518

# --- Type-changing cells (value + style via copy/paste) ---
$ws.Range("D14").Copy()
$ws.Range("G14").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("G14").PasteSpecial(-4122)

$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)

$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("D15").Copy()
$ws.Range("C29").PasteSpecial(-4163)
$ws.Range("D15").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 2

$ws.Range("D14").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)

$ws.Range("D16").Copy()
$ws.Range("C30").PasteSpecial(-4163)
$ws.Range("D16").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1

$ws.Range("C22").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("C22").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("N22").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("N22").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$excel.CutCopyMode = $false
Write-Host "applied"